$d = $word.ActiveDocument

# "Amended description of changes": the process-monitor bullet's wording
# is expanded from "...to aid filter and deletion of processes." to
# "...to aid in the filtering and deleting of processes."
$old = "Various enhancements made to the process monitor form to aid filter and deletion of processes."
$new = "Various enhancements made to the process monitor form to aid in the filtering and deleting of processes."

$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the process monitor description text to update."
}
